$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 125
$ws.Range("B125").Value2 = 6732794
$ws.Range("F125").Value2 = "FK Siauliai"
$ws.Range("G125").Value2 = "FK Dziugas Telsiai"
$ws.Range("H125").Value2 = 3
$ws.Range("K125").Value2 = 1.25
$ws.Range("L125").Value2 = 5
$ws.Range("M125").Value2 = 9
$ws.Range("N125").Value2 = 1.25
$ws.Range("O125").Value2 = 5.25
$ws.Range("P125").Value2 = 9
$ws.Range("Q125").Value2 = -1.75
$ws.Range("R125").Value2 = 2
$ws.Range("S125").Value2 = 1.8
$ws.Range("T125").Value2 = 3
$ws.Range("U125").Value2 = 1.975
$ws.Range("V125").Value2 = 1.825
$ws.Range("W125").Value2 = 0.25
$ws.Range("Z125").Value2 = 1
$ws.Range("AB125").Value2 = 0
$ws.Range("AC125").Value2 = 0

# Row 126
$ws.Range("B126").Value2 = 6732795
$ws.Range("F126").Value2 = "Suduva Marijampole"
$ws.Range("G126").Value2 = "Banga Gargzdai"
$ws.Range("H126").Value2 = 1
$ws.Range("K126").Value2 = 2.15
$ws.Range("L126").Value2 = 3.2
$ws.Range("M126").Value2 = 3
$ws.Range("N126").Value2 = 2.3
$ws.Range("O126").Value2 = 3.2
$ws.Range("P126").Value2 = 2.7
$ws.Range("Q126").Value2 = -0.25
$ws.Range("R126").Value2 = 2.05
$ws.Range("S126").Value2 = 1.75
$ws.Range("T126").Value2 = 2.25
$ws.Range("U126").Value2 = 1.9
$ws.Range("V126").Value2 = 1.9
$ws.Range("W126").Value2 = 1.3
$ws.Range("Z126").Value2 = 1.05
$ws.Range("AB126").Value2 = -1
$ws.Range("AC126").Value2 = 0.8999999999999999

# Row 164
$ws.Range("B164").Value2 = 7326568
$ws.Range("F164").Value2 = "Hegelmann Litauen"
$ws.Range("G164").Value2 = "Panevezys"
$ws.Range("I164").Value2 = 0
$ws.Range("J164").Value2 = "D"
$ws.Range("K164").Value2 = 2.375
$ws.Range("L164").Value2 = 3.2
$ws.Range("M164").Value2 = 2.625
$ws.Range("N164").Value2 = 2.7
$ws.Range("O164").Value2 = 3.2
$ws.Range("P164").Value2 = 2.3
$ws.Range("Q164").Value2 = 0
$ws.Range("R164").Value2 = 2.05
$ws.Range("S164").Value2 = 1.75
$ws.Range("T164").Value2 = 2.25
$ws.Range("U164").Value2 = 1.875
$ws.Range("V164").Value2 = 1.925
$ws.Range("X164").Value2 = 2.2
$ws.Range("Y164").Value2 = -1
$ws.Range("Z164").Value2 = 0
$ws.Range("AA164").Value2 = 0
$ws.Range("AC164").Value2 = 0.925

# Row 165
$ws.Range("B165").Value2 = 6732827
$ws.Range("F165").Value2 = "FK Dziugas Telsiai"
$ws.Range("G165").Value2 = "FK Kauno Zalgiris"
$ws.Range("I165").Value2 = 2
$ws.Range("J165").Value2 = "A"
$ws.Range("K165").Value2 = 6
$ws.Range("L165").Value2 = 3.9
$ws.Range("M165").Value2 = 1.444
$ws.Range("N165").Value2 = 4.75
$ws.Range("O165").Value2 = 3.6
$ws.Range("P165").Value2 = 1.65
$ws.Range("Q165").Value2 = 0.75
$ws.Range("R165").Value2 = 1.9
$ws.Range("S165").Value2 = 1.9
$ws.Range("T165").Value2 = 2.5
$ws.Range("U165").Value2 = 1.95
$ws.Range("V165").Value2 = 1.85
$ws.Range("X165").Value2 = -1
$ws.Range("Y165").Value2 = 0.6499999999999999
$ws.Range("Z165").Value2 = -1
$ws.Range("AA165").Value2 = 0.8999999999999999
$ws.Range("AC165").Value2 = 0.8500000000000001

# Row 177
$ws.Range("B177").Value2 = 6732837
$ws.Range("F177").Value2 = "Suduva Marijampole"
$ws.Range("G177").Value2 = "FK Riteriai"
$ws.Range("I177").Value2 = 3
$ws.Range("J177").Value2 = "A"
$ws.Range("K177").Value2 = 3.6
$ws.Range("L177").Value2 = 3.6
$ws.Range("M177").Value2 = 1.8
$ws.Range("N177").Value2 = 3
$ws.Range("O177").Value2 = 3.6
$ws.Range("P177").Value2 = 2
$ws.Range("Q177").Value2 = 0.25
$ws.Range("R177").Value2 = 2
$ws.Range("S177").Value2 = 1.8
$ws.Range("U177").Value2 = 1.975
$ws.Range("V177").Value2 = 1.825
$ws.Range("X177").Value2 = -1
$ws.Range("Y177").Value2 = 1
$ws.Range("AA177").Value2 = 0.8
$ws.Range("AB177").Value2 = 0.9750000000000001
$ws.Range("AC177").Value2 = -1

# Row 179
$ws.Range("B179").Value2 = 6732834
$ws.Range("F179").Value2 = "Panevezys"
$ws.Range("G179").Value2 = "FK Dziugas Telsiai"
$ws.Range("I179").Value2 = 0
$ws.Range("J179").Value2 = "D"
$ws.Range("K179").Value2 = 1.25
$ws.Range("L179").Value2 = 5.5
$ws.Range("M179").Value2 = 7.5
$ws.Range("N179").Value2 = 1.45
$ws.Range("O179").Value2 = 4.5
$ws.Range("P179").Value2 = 5
$ws.Range("Q179").Value2 = -1
$ws.Range("R179").Value2 = 1.775
$ws.Range("S179").Value2 = 2.025
$ws.Range("U179").Value2 = 1.875
$ws.Range("V179").Value2 = 1.925
$ws.Range("X179").Value2 = 3.5
$ws.Range("Y179").Value2 = -1
$ws.Range("AA179").Value2 = 1.025
$ws.Range("AB179").Value2 = -1
$ws.Range("AC179").Value2 = 0.925
